$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rotate weekly sample data: row2 <- old row3, row3 <- old row4, row4 <- old row2
# (Fecha, Volumen, Precio minimo/maximo/promedio, Unidad comercializacion, Precio $/Kg, Kg/unidad)

$ws.Range("D2").Value = 44855
$ws.Range("M2").Value = 25
$ws.Range("N2").Value = 15000
$ws.Range("O2").Value = 15000
$ws.Range("P2").Value = 15000
$ws.Range("Q2").Value = "$/bandeja 5 kilos"
$ws.Range("S2").Value = 3000
$ws.Range("T2").Value = 5

$ws.Range("D3").Value = 44875
$ws.Range("M3").Value = 50
$ws.Range("N3").Value = 16000
$ws.Range("O3").Value = 16000
$ws.Range("P3").Value = 16000
$ws.Range("Q3").Value = "$/bandeja 10 kilos"
$ws.Range("S3").Value = 1600
$ws.Range("T3").Value = 10

$ws.Range("D4").Value = 44874
$ws.Range("M4").Value = 67
$ws.Range("N4").Value = 16000
$ws.Range("O4").Value = 16000
$ws.Range("P4").Value = 16000
$ws.Range("Q4").Value = "$/bandeja 10 kilos"
$ws.Range("S4").Value = 1600
$ws.Range("T4").Value = 10
